# Fruta / hortaliza, semanal
# A new weekly price record is inserted at the top of the data table
# (row 21), pushing all the existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21 - this shifts rows 21:69 down to 22:70
# and mirrors the formatting of the row above (so D21 keeps the date style).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 44622
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108004
$ws.Range("J21").Value = "Papaya"
$ws.Range("K21").Value = "Cultivar IV Región"
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 75
$ws.Range("N21").Value = 25000
$ws.Range("O21").Value = 25000
$ws.Range("P21").Value = 25000
$ws.Range("Q21").Value = "$/bandeja 10 kilos"
$ws.Range("R21").Value = "Provincia del Elquí"
$ws.Range("S21").Value = 2500
$ws.Range("T21").Value = 10
